# Generate Report for Handoff
# Adds a new row for file "9e275050-23c6-49ae-919f-6d443a66c524...md" (status "Ready for handoff")
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Build the long, padded file-name strings used throughout the workbook.
# (The real filenames are a GUID padded out with repeated "o" characters to a
# fixed length, matching the pattern already used by the existing row.)
# ---------------------------------------------------------------------------
$guid = "9e275050-23c6-49ae-919f-6d443a66c524"
$mdName       = $guid.PadRight($guid.Length + 149, "o") + ".md"
$mdPath       = "e2e\" + $mdName
$zhXlfName    = $guid.PadRight($guid.Length + 40, "o") + ".f21fa609bb23d56b486f3a3405225972433bacca.zh-cn.xlf"
$deXlfName    = $guid.PadRight($guid.Length + 40, "o") + ".f21fa609bb23d56b486f3a3405225972433bacca.de-de.xlf"

$hoDateZh  = "2016-08-17 06:24:00"
$hoDateDe  = "2016-08-17 06:24:10"
$hoDateOverview = "2016-08-17 06:24:10"
$handbackEpoch = "0001-01-01 00:00:00"

$commitHash = "dd29068cfcc92195ba8d50e720e42c4ebc0c20da"
$mdTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/$mdPath"

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet: add row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value2 = $mdName
$wsOverview.Range("C3").Value2 = ".md"
$wsOverview.Range("E3").Value2 = "Ready for handoff"
$wsOverview.Range("F3").Value2 = "Ready for handoff"
$wsOverview.Range("G3").Value2 = $hoDateOverview
$wsOverview.Range("G3").NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $mdTargetUrl, "", "", $mdPath) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3")) | Out-Null

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# zh-cn sheet: add row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B3").Value2 = ".md"
$wsZh.Range("C3").Value2 = "Ready for handoff"
$wsZh.Range("D3").Value2 = "e2e"
$wsZh.Range("E3").Value2 = "ht"
$wsZh.Range("F3").Value2 = "'False"
$wsZh.Range("G3").Value2 = $zhXlfName
$wsZh.Range("H3").Value2 = $hoDateZh
$wsZh.Range("H3").NumberFormat = $dateFormat
$wsZh.Range("K3").Value2 = $handbackEpoch
$wsZh.Range("K3").NumberFormat = $dateFormat
$wsZh.Range("M3").Value2 = "'True"
$wsZh.Range("O3").Value2 = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdTargetUrl, "", "", $mdName) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3")) | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# de-de sheet: add row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B3").Value2 = ".md"
$wsDe.Range("C3").Value2 = "Ready for handoff"
$wsDe.Range("D3").Value2 = "e2e"
$wsDe.Range("E3").Value2 = "ht"
$wsDe.Range("F3").Value2 = "'False"
$wsDe.Range("G3").Value2 = $deXlfName
$wsDe.Range("H3").Value2 = $hoDateDe
$wsDe.Range("H3").NumberFormat = $dateFormat
$wsDe.Range("K3").Value2 = $handbackEpoch
$wsDe.Range("K3").NumberFormat = $dateFormat
$wsDe.Range("M3").Value2 = "'True"
$wsDe.Range("O3").Value2 = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdTargetUrl, "", "", $mdName) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3")) | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = 16.3
